# "Added BarHomogeneous to gurobi solver" -- the solver change altered the
# optimisation results, which were then re-exported into this summary
# workbook. Apply the updated values and drop the now-unused
# "Total non-actualized Operation cost [USD]" row from the first sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Project Total Costs" ---------------------------------------
$ws1 = $wb.Worksheets.Item("Project Total Costs")

# Remove the "Total non-actualized Operation cost [USD]" row entirely.
$ws1.Rows.Item(4).Delete()

$ws1.Range("B2").Value = 373723.33889900002
$ws1.Range("B3").Value = 322513.12391899998
$ws1.Range("B4").Value = 384023.54050200002
$ws1.Range("B5").Value = 332813.32552299998
$ws1.Range("B6").Value = 0.3571982167907834

# --- Sheet "Components Capacity and Cost" -------------------------------
$ws2 = $wb.Worksheets.Item("Components Capacity and Cost")

$ws2.Range("B2").Value = 0.00053158528266100008
$ws2.Range("B3").Value = 70.604926507499997
$ws2.Range("B4").Value = 43.435200063400004
$ws2.Range("B5").Value = 92.409090404799997
$ws2.Range("B6").Value = 0.000000055650272516400002
$ws2.Range("B7").Value = 0.39337310916913998
$ws2.Range("B8").Value = 29661.12962580075
$ws2.Range("B9").Value = 13030.56001902
$ws2.Range("B10").Value = 341331.45722820982
$ws2.Range("B11").Value = 0.00025573247980823871
$ws2.Range("B12").Value = 384023.54050187219

# --- Sheet "Yearly Costs Info" ------------------------------------------
$ws3 = $wb.Worksheets.Item("Yearly Costs Info")

$ws3.Range("B2").Value = 0.0078674621833828007
$ws3.Range("C2").Value = 4138.8633642918749
$ws3.Range("D2").Value = 6826.6291496788454
$ws3.Range("E2").Value = 10965.500381432899
$ws3.Range("F2").Value = 150291.0087386785
$ws3.Range("G2").Value = 0.056681316501663033

$ws3.Range("B3").Value = 0.0078674621833828007
$ws3.Range("C3").Value = 4138.8633642918749
$ws3.Range("D3").Value = 6826.6291496788454
$ws3.Range("E3").Value = 10965.500381432899
$ws3.Range("F3").Value = 150291.0010620206
$ws3.Range("G3").Value = 0.056674539760673173

# --- Sheet "Yearly Energy Averages" --------------------------------------
$ws4 = $wb.Worksheets.Item("Yearly Energy Averages")

$ws4.Range("B2").Value = 63.042880449778529
$ws4.Range("C2").Value = 5.0572701690140596
$ws4.Range("D2").Value = 31.899855489291699
$ws4.Range("E2").Value = 10.069713626069801

$ws4.Range("B3").Value = 63.042938341288803
$ws4.Range("C3").Value = 5.0572121527303668
$ws4.Range("D3").Value = 31.899855609597719
$ws4.Range("E3").Value = 10.069713529177401
